# This edit swaps the data content of rows 5 and 6 on the "Artfynd" sheet.
# (Row 5 becomes what row 6 used to be, and vice versa, for the columns
#  that actually differ between the two records.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the current ("before") values of the columns that differ
# between row 5 and row 6 so we can swap them safely.
# (Value2 is used for reading because it reliably returns the raw
#  scalar data in this environment.)
$A5 = $ws.Range("A5").Value2
$A6 = $ws.Range("A6").Value2

$B5 = $ws.Range("B5").Value2
$B6 = $ws.Range("B6").Value2

$E5 = $ws.Range("E5").Value2
$E6 = $ws.Range("E6").Value2

$F5 = $ws.Range("F5").Value2
$F6 = $ws.Range("F6").Value2

$G5 = $ws.Range("G5").Value2
$G6 = $ws.Range("G6").Value2

$H5 = $ws.Range("H5").Value2
$H6 = $ws.Range("H6").Value2

$M5 = $ws.Range("M5").Value2
$M6 = $ws.Range("M6").Value2

$Q5 = $ws.Range("Q5").Value2
$Q6 = $ws.Range("Q6").Value2

$R5 = $ws.Range("R5").Value2
$R6 = $ws.Range("R6").Value2

# Write row 5 with the former row 6 values.
$ws.Range("A5").Value = $A6
$ws.Range("B5").Value = $B6
$ws.Range("E5").Value = $E6
$ws.Range("F5").Value = $F6
$ws.Range("G5").Value = $G6
$ws.Range("H5").Value = $H6
$ws.Range("M5").Value = $M6
$ws.Range("Q5").Value = $Q6
$ws.Range("R5").Value = $R6

# Write row 6 with the former row 5 values.
$ws.Range("A6").Value = $A5
$ws.Range("B6").Value = $B5
$ws.Range("E6").Value = $E5
$ws.Range("F6").Value = $F5
$ws.Range("G6").Value = $G5
$ws.Range("H6").Value = $H5
$ws.Range("M6").Value = $M5
$ws.Range("Q6").Value = $Q5
$ws.Range("R6").Value = $R5
